# Adjust input endTime to capture the whole day if time is not entered.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 1900
$ws.Range("C20").Value = 190

$ws.Range("B21").Value = 20000
$ws.Range("C21").Value = 200

$ws.Range("B22").Value = 21000
$ws.Range("C22").Value = 210

$ws.Range("C23").Select()
